# Add "first class soul" localization entries (SOUL_NAME_*/SOUL_DESC_*) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 30,2
$arr[0,0] = "SOUL_NAME_1000"
$arr[0,1] = "狗"
$arr[1,0] = "SOUL_DESC_1000"
$arr[1,1] = "狗"
$arr[2,0] = "SOUL_NAME_1001"
$arr[2,1] = "猫"
$arr[3,0] = "SOUL_DESC_1001"
$arr[3,1] = "猫"
$arr[4,0] = "SOUL_NAME_1002"
$arr[4,1] = "熊"
$arr[5,0] = "SOUL_DESC_1002"
$arr[5,1] = "熊"
$arr[6,0] = "SOUL_NAME_2001"
$arr[6,1] = "拳师"
$arr[7,0] = "SOUL_DESC_2001"
$arr[7,1] = "拳师"
$arr[8,0] = "SOUL_NAME_2002"
$arr[8,1] = "武士"
$arr[9,0] = "SOUL_DESC_2002"
$arr[9,1] = "武士"
$arr[10,0] = "SOUL_NAME_2003"
$arr[10,1] = "学者"
$arr[11,0] = "SOUL_DESC_2003"
$arr[11,1] = "学者"
$arr[12,0] = "SOUL_NAME_3000"
$arr[12,1] = "兽人"
$arr[13,0] = "SOUL_DESC_3000"
$arr[13,1] = "兽人"
$arr[14,0] = "SOUL_NAME_3001"
$arr[14,1] = "牛头人"
$arr[15,0] = "SOUL_DESC_3001"
$arr[15,1] = "牛头人"
$arr[16,0] = "SOUL_NAME_3002"
$arr[16,1] = "兔头萨满"
$arr[17,0] = "SOUL_DESC_3002"
$arr[17,1] = "兔头萨满"
$arr[18,0] = "SOUL_NAME_4000"
$arr[18,1] = "精灵"
$arr[19,0] = "SOUL_DESC_4000"
$arr[19,1] = "精灵"
$arr[20,0] = "SOUL_NAME_4001"
$arr[20,1] = "古树"
$arr[21,0] = "SOUL_DESC_4001"
$arr[21,1] = "古树"
$arr[22,0] = "SOUL_NAME_4002"
$arr[22,1] = "黑精灵"
$arr[23,0] = "SOUL_DESC_4002"
$arr[23,1] = "黑精灵"
$arr[24,0] = "SOUL_NAME_5000"
$arr[24,1] = "骷髅"
$arr[25,0] = "SOUL_DESC_5000"
$arr[25,1] = "骷髅"
$arr[26,0] = "SOUL_NAME_5001"
$arr[26,1] = "僵尸"
$arr[27,0] = "SOUL_DESC_5001"
$arr[27,1] = "僵尸"
$arr[28,0] = "SOUL_NAME_5002"
$arr[28,1] = "鬼魂"
$arr[29,0] = "SOUL_DESC_5002"
$arr[29,1] = "鬼魂"
$ws.Range("A63:B92").Value = $arr

# Match the author's final selection / scroll position.
$ws.Range("D88").Select()
